$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Worksheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 89
$ws1.Range("F5").Value = 30
$ws1.Range("F8").Value = 8049
$ws1.Range("F9").Value = 759
$ws1.Range("F10").Value = 241
$ws1.Range("F12").Value = 786
$ws1.Range("F13").Value = 35
$ws1.Range("F19").Value = 856

# --- Sheet "全部类型" (Worksheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 89
$ws4.Range("F5").Value = 30
$ws4.Range("F9").Value = 8049
$ws4.Range("F10").Value = 759
$ws4.Range("F11").Value = 241
$ws4.Range("F13").Value = 786
$ws4.Range("F14").Value = 35
$ws4.Range("F20").Value = 856
